# Add new event rows (228-234) to the events worksheet, matching the
# data/formatting of the other filled rows in the sheet.
#
# Note: this COM-interop runtime's Value assignment is only recognised at
# top-level statement scope, not inside a custom `function`/`param()`
# block, so the per-row logic below is a flat `foreach` loop over a data
# table rather than a helper function.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$events = @(
    @{ Row = 228; Date = 45738; Event = "BYORN"; Location = "Junkyard"; City = "Dortmund"; Link = "https://www.instagram.com/reel/DF5on0etx8D/?igsh=MTNocWxkcDJtemdqNQ==" },
    @{ Row = 229; Date = 45808; Event = "ONEYEAR PRISMATICCLUB 12h RAVE"; Location = "Prismatic"; City = "Dortmund"; Link = "https://www.instagram.com/reel/DG-bQ4gts2V/?igsh=MW5nM2RvaGFyaHlzMA==" },
    @{ Row = 230; Date = 45729; Event = "180 MIN RAVE (20-23Uhr)"; Location = "Prismatic"; City = "Dortmund"; Link = "https://www.instagram.com/reel/DG6QYyPNsrU/?igsh=MTRzY2I3aHFuODVnZQ==" },
    @{ Row = 231; Date = 45731; Event = "PRESENTED BY VOIT: DEADLY SINS"; Location = "SNRS Club"; City = "Dortmund"; Link = "https://www.instagram.com/reel/DGizSy8qhqO/?igsh=bHowMXlmaXVkZzlw" },
    @{ Row = 232; Date = 45751; Event = "EHRENKLUB"; Location = "Odonien"; City = "Köln"; Link = "https://www.instagram.com/p/DG6BJ_VMhkd/?igsh=MXA5Nzl1bnpsbjF2Zg==" },
    @{ Row = 233; Date = 45767; Event = "EHRENKLUB OSTERSONNTAG"; Location = "Schrotty"; City = "Köln"; Link = "https://www.instagram.com/p/DFx-VlKgWs-/?igsh=bGl0N2R1ejIyYXJv" },
    @{ Row = 234; Date = 45744; Event = "DYSPHORIA x SYNCED"; Location = "Rotunde"; City = "Bochum"; Link = "https://www.instagram.com/reel/DGvwBc1NMZx/?igsh=MXFoMWQ0dmhqbnk1dw==" }
)

foreach ($e in $events) {
    $r = $e.Row

    # Copy the formatting (number format / fill / border / font) from the
    # last filled-in row (227) onto the B:E cells of this row first, so the
    # new cells pick up the same style as the other data rows instead of
    # keeping the blank placeholder style.
    $ws.Range("B227:E227").Copy() | Out-Null
    $ws.Range("B" + $r + ":E" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $r).Value = $e.Date
    $ws.Range("B" + $r).Value = $e.Event
    $ws.Range("C" + $r).Value = $e.Location
    $ws.Range("D" + $r).Value = $e.City

    $ws.Hyperlinks.Add($ws.Range("E" + $r), $e.Link, "", "", $e.Link) | Out-Null

    # Hyperlinks.Add re-styles the target cell with the default hyperlink
    # look (blue/underline). Restore the plain data-row style by
    # re-pasting formats from a plain (non-hyperlink) cell of row 227.
    $ws.Range("D227").Copy() | Out-Null
    $ws.Range("E" + $r).PasteSpecial(-4122) | Out-Null
}

$ws.Range("A1").Select() | Out-Null

Write-Host "Added rows 228-234"
